$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(227, 44301, 4, 36, 210.0962941348118),
    @(228, 44302, 9, 39, 227.6043186460461),
    @(229, 44303, 7, 43, 250.9483513276919)
)

foreach ($row in $newRows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]

    # Column A in the data rows above uses the bold/bordered/centered date style (s="2");
    # copy that formatting down onto the newly appended rows.
    $ws.Cells.Item($r - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
}
